$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $rng = $Sheet.Range($Addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $Val
    $rng.Style = $origStyle
}

Set-TextValue $ws 'D2' '26.959.82'
Set-TextValue $ws 'E2' '  +0.48%  '
Set-TextValue $ws 'D3' '1.556.53'
Set-TextValue $ws 'E3' '  -0.28%  '
Set-TextValue $ws 'E4' '  +0.48%  '
Set-TextValue $ws 'D5' '207.02'
Set-TextValue $ws 'E5' '  +0.75%  '
Set-TextValue $ws 'E6' '  +1.18%  '
Set-TextValue $ws 'E7' '  +0.43%  '
Set-TextValue $ws 'E8' '  +0.67%  '
Set-TextValue $ws 'D9' '21.54'
Set-TextValue $ws 'E9' '  -0.09%  '
Set-TextValue $ws 'E10' '  -0.25%  '
Set-TextValue $ws 'D11' '0.0859'
Set-TextValue $ws 'E11' '  -0.22%  '
Set-TextValue $ws 'D12' '1.778.74'
Set-TextValue $ws 'E12' '  -0.17%  '
Set-TextValue $ws 'D13' '1.559.26'
Set-TextValue $ws 'E13' '  -0.75%  '
Set-TextValue $ws 'D14' '3.70'
Set-TextValue $ws 'E14' '  -0.43%  '
Set-TextValue $ws 'E15' '  +0.27%  '
Set-TextValue $ws 'D16' '26.961.27'
Set-TextValue $ws 'E16' '  +0.47%  '
Set-TextValue $ws 'D17' '61.77'
Set-TextValue $ws 'E17' '  +0.82%  '
Set-TextValue $ws 'D18' '214.75'
Set-TextValue $ws 'E18' '  -0.19%  '
Set-TextValue $ws 'D19' '0.0₃0687'
Set-TextValue $ws 'E19' '  +0.86%  '
Set-TextValue $ws 'E20' '  -1.18%  '
Set-TextValue $ws 'E21' '  +0.42%  '
Set-TextValue $ws 'E22' '  -1.84%  '
Set-TextValue $ws 'D23' '9.21'
Set-TextValue $ws 'E23' '  +0.58%  '
Set-TextValue $ws 'E24' '  -2.13%  '
Set-TextValue $ws 'D25' '153.39'
Set-TextValue $ws 'E25' '  +0.02%  '
Set-TextValue $ws 'D26' '6.66'
Set-TextValue $ws 'E26' '  +0.75%  '
Set-TextValue $ws 'D27' '14.89'
Set-TextValue $ws 'E27' '  -0.98%  '
Set-TextValue $ws 'E28' '  +0.40%  '
Set-TextValue $ws 'E29' '  +0.52%  '
Set-TextValue $ws 'D30' '0.0459'
Set-TextValue $ws 'E30' '  -1.32%  '
Set-TextValue $ws 'E31' '  -0.44%  '
Set-TextValue $ws 'E32' '  +1.87%  '
Set-TextValue $ws 'D33' '1.374.33'
Set-TextValue $ws 'E33' '  -0.53%  '
Set-TextValue $ws 'E34' '  +1.48%  '
Set-TextValue $ws 'E35' '  +2.37%  '
Set-TextValue $ws 'D36' '0.969'
Set-TextValue $ws 'E36' '  +5.64%  '
Set-TextValue $ws 'E38' '  +0.84%  '
Set-TextValue $ws 'D39' '0.520'
Set-TextValue $ws 'E39' '  -1.25%  '
Set-TextValue $ws 'D40' '0.808'
Set-TextValue $ws 'E40' '  -0.14%  '
Set-TextValue $ws 'E41' '  +0.42%  '
Set-TextValue $ws 'D42' '5.52'
Set-TextValue $ws 'E42' '  -0.02%  '
Set-TextValue $ws 'E43' '  -1.26%  '
Set-TextValue $ws 'D44' '2.24'
Set-TextValue $ws 'E44' '  +2.65%  '
Set-TextValue $ws 'D45' '63.76'
Set-TextValue $ws 'E45' '  +0.43%  '
Set-TextValue $ws 'D46' '1.74'
Set-TextValue $ws 'E46' '  -1.81%  '
Set-TextValue $ws 'D47' '1.691.94'
Set-TextValue $ws 'E47' '  -0.26%  '
Set-TextValue $ws 'E48' '  -3.46%  '
Set-TextValue $ws 'D49' '86.33'
Set-TextValue $ws 'E49' '  -0.20%  '
Set-TextValue $ws 'E50' '  +0.05%  '
Set-TextValue $ws 'B51' 'Algorand'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D51' '0.0956'
Set-TextValue $ws 'E51' '  +0.50%  '
